# Apply the "想去人数"/"最低票价" data refresh described in the commit
# "Update gh-pages to output generated at 456a3b4": numeric counters bump
# up on sheets 展览 (Exhibition), 演出 (Performance), 本地生活 (Local life)
# and 全部类型 (All types), and two previously-unavailable-ticket cells
# ("不可售") now have a numeric price.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet 1) ---
$ws1.Range("F5").Value = 2497
$ws1.Range("G6").Value = 158
$ws1.Range("F8").Value = 3175
$ws1.Range("F10").Value = 4668
$ws1.Range("F14").Value = 600
$ws1.Range("F17").Value = 5
$ws1.Range("F20").Value = 15
$ws1.Range("F24").Value = 4651
$ws1.Range("F25").Value = 11
$ws1.Range("F28").Value = 5492
$ws1.Range("F29").Value = 11
$ws1.Range("F30").Value = 1168
$ws1.Range("F32").Value = 640
$ws1.Range("F33").Value = 4398
$ws1.Range("F35").Value = 65
$ws1.Range("F37").Value = 762
$ws1.Range("F39").Value = 699
$ws1.Range("F40").Value = 701

# --- 演出 (sheet 2) ---
$ws2.Range("F6").Value = 42

# --- 本地生活 (sheet 3) ---
$ws3.Range("F4").Value = 30

# --- 全部类型 (sheet 4) ---
$ws4.Range("F5").Value = 30
$ws4.Range("F8").Value = 2497
$ws4.Range("G9").Value = 158
$ws4.Range("F12").Value = 3175
$ws4.Range("F14").Value = 4668
$ws4.Range("F18").Value = 600
$ws4.Range("F21").Value = 5
$ws4.Range("F24").Value = 15
$ws4.Range("F29").Value = 4651
$ws4.Range("F30").Value = 11
$ws4.Range("F33").Value = 5492
$ws4.Range("F34").Value = 11
$ws4.Range("F35").Value = 1168
$ws4.Range("F37").Value = 640
$ws4.Range("F38").Value = 4398
$ws4.Range("F41").Value = 65
$ws4.Range("F43").Value = 762
$ws4.Range("F45").Value = 699
$ws4.Range("F46").Value = 701
$ws4.Range("F48").Value = 42
